# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Updates the "K" column (column G) values for rows 2-37 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 6
    3  = 6
    4  = 3
    5  = 5
    6  = 11
    7  = 6
    8  = 4
    9  = 5
    10 = 7
    11 = 6
    12 = 9
    13 = 5
    14 = 4
    15 = 7
    16 = 5
    17 = 9
    18 = 9
    19 = 6
    20 = 7
    21 = 5
    22 = 5
    23 = 8
    24 = 6
    25 = 4
    26 = 4
    27 = 7
    28 = 2
    29 = 8
    30 = 4
    31 = 6
    32 = 3
    33 = 3
    34 = 8
    35 = 5
    36 = 5
    37 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
